# The author changed the table style applied to the financial-documents
# table on slide 5 (the 3-column table, shape index 2 on that slide) from
# the locally-defined "Table_0" style to a different (built-in) table
# style, identified by its style GUID.
#
# Table styles in PowerPoint's object model are not assigned by setting a
# writable property - attempting `$tbl.Style = "{GUID}"` raises
# "Table styles cannot be assigned through a property - call
# Table.ApplyStyle("{GUID}") instead", so we use Table.ApplyStyle, which is
# exactly how PowerPoint itself re-writes <a:tableStyleId> under
# <a:tblPr> in the slide's graphicFrame/table XML.

$p = $ppt.ActivePresentation

# Slide 5, 2nd shape ("Google Shape;122;p17") is the graphicFrame holding
# the 3x? table whose <a:tableStyleId> changes in the diff.
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)

$table = $tableShape.Table
$table.ApplyStyle("{EDBCFD0D-97B2-4782-B75D-9E13DC3CB70A}")
